$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction in SA algorithm fitness log values (run_14, mp0-8_5)
$ws.Range("C2:C5").Value = 9184
$ws.Range("C6:C18").Value = 8541
$ws.Range("C19").Value = 8302
$ws.Range("C20:C51").Value = 7667
$ws.Range("C52:C57").Value = 7320
$ws.Range("C58:C129").Value = 7310
$ws.Range("C190:C242").Value = 7293

Write-Output "Updated fitness values"